$p = $ppt.ActivePresentation

# --- Slide 3: "Input graph into LLM" bullet becomes bold (minus the "**" markers) ---
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange

$para1 = $tr3.Paragraphs(1, 1)
$chars1 = $para1.Characters(1, $para1.Text.Length)
$chars1.Text = "Input graph into LLM:"
$boldPart1 = $para1.Characters(1, 20)
$boldPart1.Font.Bold = $true

# --- Slide 4: "Applications - Pure graphs", Build-a-Graph / Context-Summarization wording, and closing line ---
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange

# Paragraph 1: "** Applications - Pure graphs **" -> bold "Applications - Pure graphs", plus extra space-before
$paraApps = $tr4.Paragraphs(1, 1)
$paraApps.ParagraphFormat.SpaceBefore = 30
$paraApps.ParagraphFormat.Bullet.Type = 1
$paraApps.ParagraphFormat.Bullet.Type = 0
$charsApps = $paraApps.Characters(1, $paraApps.Text.Length)
$charsApps.Text = "Applications - Pure graphs"
$paraApps.Font.Bold = $true

# Paragraph 10: Build-a-Graph wording tweak
$paraBuild = $tr4.Paragraphs(10, 1)
$charsBuild = $paraBuild.Characters(1, $paraBuild.Text.Length)
$charsBuild.Text = "Build-a-Graph: reconstruct the relevant graph structures and then perform reasoning on them. This"

# Paragraph 11: Context-Summarization wording tweak
$paraCtx = $tr4.Paragraphs(11, 1)
$charsCtx = $paraCtx.Characters(1, $paraCtx.Text.Length)
$charsCtx.Text = "Context-Summarization: summarize the key nodes, edges, or sub-graphs and perform reasoning."

# Paragraph 15: "** Overall, no consensus on how to represent graphs **" -> bold, no "**"
$paraOverall = $tr4.Paragraphs(15, 1)
$charsOverall = $paraOverall.Characters(1, $paraOverall.Text.Length)
$charsOverall.Text = "Overall, no consensus on how to represent graphs"
$paraOverall.Font.Bold = $true
